# Daily attendance processing - 2025-11-11 04:50:49
# Reorders the "Recorded By" (column G) entries so that entries of the
# form "System, <email>" become "<email>, System" (System is moved to the
# end of the list). Entries that already have a different order, a single
# name, three or more names, or specifically pair "System" with
# "backup@backdoor.com" are left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    if ($val.StartsWith("System, ")) {
        $rest = $val.Substring(8)

        # Keep "System, backup@backdoor.com" unchanged, and only swap
        # simple two-party "System, X" pairs (no further commas in $rest).
        if ($rest -ne "backup@backdoor.com" -and -not ($rest.Contains(","))) {
            $cell.Value = $rest + ", System"
        }
    }
}
